$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the two new "linked maps" sheets at the right tab positions:
#    bray_baseline, bray_lm, louisville_baseline, louisville_lm
# ---------------------------------------------------------------------------
$brayBaseline = $wb.Worksheets.Item("bray_baseline")
$brayLm = $wb.Worksheets.Add([System.Type]::Missing, $brayBaseline)
$brayLm.Name = "bray_lm"

$louisvilleBaseline = $wb.Worksheets.Item("louisville_baseline")
$louisvilleLm = $wb.Worksheets.Add([System.Type]::Missing, $louisvilleBaseline)
$louisvilleLm.Name = "louisville_lm"

# ---------------------------------------------------------------------------
# 2. Column widths / layout matching the baseline sheets
# ---------------------------------------------------------------------------
$brayLm.Columns.Item(1).ColumnWidth = 25.5
$brayLm.Columns.Item(2).ColumnWidth = 39.33203125
$brayLm.Columns.Item(3).ColumnWidth = 32.6640625
$brayLm.Columns.Item(4).ColumnWidth = 36.33203125

$louisvilleLm.Columns.Item(1).ColumnWidth = 25.5
$louisvilleLm.Columns.Item(2).ColumnWidth = 39.33203125
$louisvilleLm.Columns.Item(3).ColumnWidth = 32.6640625
$louisvilleLm.Columns.Item(4).ColumnWidth = 36.33203125

# ---------------------------------------------------------------------------
# 3. Header row (same headers/styles as the other sheets)
# ---------------------------------------------------------------------------
function Set-HeaderRow($ws) {
    $ws.Range("A1").Value = "GID"
    $ws.Range("B1").Value = "Instance"
    $ws.Range("C1").Value = "Label"
    $ws.Range("D1").Value = "Types"
    $ws.Range("E1").Value = "Match"
    $ws.Range("A1:C1").Style = "Heading 1"
    $ws.Range("A1:E1").Font.Bold = $true
    $ws.Range("A1:E1").HorizontalAlignment = -4108
    $ws.Range("A1:E1").VerticalAlignment = -4160
    $ws.Range("A1:E1").Borders.LineStyle = 1
}

Set-HeaderRow $brayLm
Set-HeaderRow $louisvilleLm

# ---------------------------------------------------------------------------
# 4. bray_lm data rows
# ---------------------------------------------------------------------------
$brayLmRows = @(
    @("http://linkedmaps.isi.edu/75", "http://linkedgeodata.org/triplify/way10670872", "Black Butte Subdivision", "http://geovocab.org/spatial#Feature http://linkedgeodata.org/meta/Way http://linkedgeodata.org/ontology/Rail http://linkedgeodata.org/ontology/RailwayThing", 1),
    @("http://linkedmaps.isi.edu/75", "http://linkedgeodata.org/triplify/way249503576", "Black Butte Subdivision", "http://geovocab.org/spatial#Feature http://linkedgeodata.org/meta/Way http://linkedgeodata.org/ontology/Rail http://linkedgeodata.org/ontology/RailwayThing", 1),
    @("http://linkedmaps.isi.edu/75", "http://linkedgeodata.org/triplify/way322131253", "Black Butte Subdivision", "http://geovocab.org/spatial#Feature http://linkedgeodata.org/meta/Way http://linkedgeodata.org/ontology/Rail http://linkedgeodata.org/ontology/RailwayThing", 1),
    @("http://linkedmaps.isi.edu/75", "http://linkedgeodata.org/triplify/way249503577", "Black Butte Subdivision", "http://geovocab.org/spatial#Feature http://linkedgeodata.org/meta/Way http://linkedgeodata.org/ontology/Rail http://linkedgeodata.org/ontology/RailwayThing", 1),
    @("http://linkedmaps.isi.edu/81", "http://linkedgeodata.org/triplify/way322131262", "Black Butte Subdivision", "http://geovocab.org/spatial#Feature http://linkedgeodata.org/meta/Way http://linkedgeodata.org/ontology/Rail http://linkedgeodata.org/ontology/RailwayThing", 1),
    @("http://linkedmaps.isi.edu/81", "http://linkedgeodata.org/triplify/way249503577", "Black Butte Subdivision", "http://geovocab.org/spatial#Feature http://linkedgeodata.org/meta/Way http://linkedgeodata.org/ontology/Rail http://linkedgeodata.org/ontology/RailwayThing", 0),
    @("http://linkedmaps.isi.edu/81", "http://linkedgeodata.org/triplify/way10670906", "Black Butte Subdivision", "http://geovocab.org/spatial#Feature http://linkedgeodata.org/meta/Way http://linkedgeodata.org/ontology/Rail http://linkedgeodata.org/ontology/RailwayThing", 1),
    @("http://linkedmaps.isi.edu/80", "http://linkedgeodata.org/triplify/way10670920", "Black Butte Subdivision", "http://geovocab.org/spatial#Feature http://linkedgeodata.org/meta/Way http://linkedgeodata.org/ontology/Rail http://linkedgeodata.org/ontology/RailwayThing", 1),
    @("http://linkedmaps.isi.edu/80", "http://linkedgeodata.org/triplify/way322131253", "Black Butte Subdivision", "http://geovocab.org/spatial#Feature http://linkedgeodata.org/meta/Way http://linkedgeodata.org/ontology/Rail http://linkedgeodata.org/ontology/RailwayThing", 0),
    @("http://linkedmaps.isi.edu/69", "http://linkedgeodata.org/triplify/way177559138", "Long Bell Lumber Company Railroad", "http://geovocab.org/spatial#Feature http://linkedgeodata.org/meta/Way http://linkedgeodata.org/ontology/AbandonedRailway http://linkedgeodata.org/ontology/RailwayThing", 0),
    @("http://linkedmaps.isi.edu/69", "http://linkedgeodata.org/triplify/way177559134", "Long Bell Lumber Company Railroad", "http://geovocab.org/spatial#Feature http://linkedgeodata.org/meta/Way http://linkedgeodata.org/ontology/AbandonedRailway http://linkedgeodata.org/ontology/RailwayThing", 1),
    @("http://linkedmaps.isi.edu/71", "http://linkedgeodata.org/triplify/way177559138", "Long Bell Lumber Company Railroad", "http://geovocab.org/spatial#Feature http://linkedgeodata.org/meta/Way http://linkedgeodata.org/ontology/AbandonedRailway http://linkedgeodata.org/ontology/RailwayThing", 1),
    @("http://linkedmaps.isi.edu/71", "http://linkedgeodata.org/triplify/way10661139", "Long Bell Lumber Company Railroad", "http://geovocab.org/spatial#Feature http://linkedgeodata.org/meta/Way http://linkedgeodata.org/ontology/AbandonedRailway http://linkedgeodata.org/ontology/RailwayThing", 1),
    @("http://linkedmaps.isi.edu/83", "http://linkedgeodata.org/triplify/way249503576", "Black Butte Subdivision", "http://geovocab.org/spatial#Feature http://linkedgeodata.org/meta/Way http://linkedgeodata.org/ontology/Rail http://linkedgeodata.org/ontology/RailwayThing", 0),
    @("http://linkedmaps.isi.edu/79", "http://linkedgeodata.org/triplify/way10661144", "Long Bell Lumber Company Railroad", "http://geovocab.org/spatial#Feature http://linkedgeodata.org/meta/Way http://linkedgeodata.org/ontology/Rail http://linkedgeodata.org/ontology/RailwayThing", 1)
)

$r = 2
foreach ($row in $brayLmRows) {
    $brayLm.Cells.Item($r, 1).Value = $row[0]
    $brayLm.Cells.Item($r, 2).Value = $row[1]
    $brayLm.Cells.Item($r, 3).Value = $row[2]
    $brayLm.Cells.Item($r, 4).Value = $row[3]
    $brayLm.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# Rows 13 & 14 (instance column) render bold in the source workbook
$brayLm.Range("B13").Font.Bold = $true
$brayLm.Range("B14").Font.Bold = $true

# ---------------------------------------------------------------------------
# 5. louisville_lm data rows
# ---------------------------------------------------------------------------
$louisvilleLmRows = @(
    @("http://linkedmaps.isi.edu/4", "http://linkedgeodata.org/triplify/way17071151", "Rocky Flats Industrial Lead", "http://geovocab.org/spatial#Feature http://linkedgeodata.org/meta/Way http://linkedgeodata.org/ontology/Rail http://linkedgeodata.org/ontology/RailwayThing", 1),
    @("http://linkedmaps.isi.edu/4", "http://linkedgeodata.org/triplify/way39723232", "Rocky Flats Industrial Lead", "http://geovocab.org/spatial#Feature http://linkedgeodata.org/meta/Way http://linkedgeodata.org/ontology/Rail http://linkedgeodata.org/ontology/RailwayThing", 1),
    @("http://linkedmaps.isi.edu/4", "http://linkedgeodata.org/triplify/way39723233", "Rocky Flats Industrial Lead", "http://geovocab.org/spatial#Feature http://linkedgeodata.org/meta/Way http://linkedgeodata.org/ontology/Rail http://linkedgeodata.org/ontology/RailwayThing", 1),
    @("http://linkedmaps.isi.edu/4", "http://linkedgeodata.org/triplify/way17071155", "Rocky Flats Industrial Lead", "http://geovocab.org/spatial#Feature http://linkedgeodata.org/meta/Way http://linkedgeodata.org/ontology/Rail http://linkedgeodata.org/ontology/RailwayThing", 0),
    @("http://linkedmaps.isi.edu/10", "http://linkedgeodata.org/triplify/way43007631", "BNSF Front Range Subdivision", "http://geovocab.org/spatial#Feature http://linkedgeodata.org/meta/Way http://linkedgeodata.org/ontology/Rail http://linkedgeodata.org/ontology/RailwayThing", 1),
    @("http://linkedmaps.isi.edu/10", "http://linkedgeodata.org/triplify/way43007632", "Front Range Subdivision", "http://geovocab.org/spatial#Feature http://linkedgeodata.org/meta/Way http://linkedgeodata.org/ontology/Rail http://linkedgeodata.org/ontology/RailwayThing", 1)
)

$r = 2
foreach ($row in $louisvilleLmRows) {
    $louisvilleLm.Cells.Item($r, 1).Value = $row[0]
    $louisvilleLm.Cells.Item($r, 2).Value = $row[1]
    $louisvilleLm.Cells.Item($r, 3).Value = $row[2]
    $louisvilleLm.Cells.Item($r, 4).Value = $row[3]
    $louisvilleLm.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 6. Selections / active sheet, matching the saved cursor state in the diff
# ---------------------------------------------------------------------------
$brayLm.Range("A17").Select()
$louisvilleLm.Range("A2").Select()

$louisvilleLm.Activate()
